$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/centered/bordered
# style (style index 1) used by every other header cell in row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the team record (Wins/Losses/Ties) for every data row (2-51).
$ws.Range("AD2:AD51").Value = 79
$ws.Range("AE2:AE51").Value = 83
$ws.Range("AF2:AF51").Value = 0
